{"js": "// Replace the three-digit \u00f7 one-digit division answers in each table cell\n// with the new set of problems/answers, per the commit diff.\nconst replacements = [\n  [\"572\u00f72=286, 0\", \"375\u00f78=46, 7\"],\n  [\"890\u00f76=148, 2\", \"621\u00f72=310, 1\"],\n  [\"479\u00f75=95, 4\", \"973\u00f73=324, 1\"],\n  [\"758\u00f72=379, 0\", \"311\u00f79=34, 5\"],\n  [\"261\u00f75=52, 1\", \"939\u00f75=187, 4\"],\n  [\"577\u00f74=144, 1\", \"573\u00f75=114, 3\"],\n  [\"903\u00f79=100, 3\", \"564\u00f78=70, 4\"],\n  [\"545\u00f74=136, 1\", \"797\u00f75=159, 2\"],\n  [\"576\u00f78=72, 0\", \"954\u00f73=318, 0\"],\n  [\"744\u00f72=372, 0\", \"230\u00f76=38, 2\"],\n  [\"793\u00f76=132, 1\", \"148\u00f76=24, 4\"],\n  [\"528\u00f78=66, 0\", \"290\u00f75=58, 0\"],\n  [\"108\u00f79=12, 0\", \"340\u00f74=85, 0\"],\n  [\"749\u00f76=124, 5\", \"913\u00f76=152, 1\"],\n  [\"224\u00f74=56, 0\", \"658\u00f78=82, 2\"],\n  [\"260\u00f78=32, 4\", \"925\u00f75=185, 0\"],\n  [\"607\u00f73=202, 1\", \"286\u00f76=47, 4\"],\n  [\"927\u00f78=115, 7\", \"267\u00f72=133, 1\"],\n  [\"923\u00f78=115, 3\", \"902\u00f79=100, 2\"],\n  [\"350\u00f72=175, 0\", \"380\u00f78=47, 4\"],\n  [\"121\u00f77=17, 2\", \"807\u00f79=89, 6\"],\n  [\"219\u00f73=73, 0\", \"294\u00f73=98, 0\"],\n  [\"991\u00f75=198, 1\", \"706\u00f72=353, 0\"],\n  [\"683\u00f75=136, 3\", \"980\u00f73=326, 2\"],\n  [\"254\u00f74=63, 2\", \"732\u00f78=91, 4\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const found = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  if (found.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const range of found.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# Old => new three-digit / one-digit division problems, in document order.\n$pairs = @(\n    @(\"572\u00f72=286, 0\", \"375\u00f78=46, 7\"),\n    @(\"890\u00f76=148, 2\", \"621\u00f72=310, 1\"),\n    @(\"479\u00f75=95, 4\", \"973\u00f73=324, 1\"),\n    @(\"758\u00f72=379, 0\", \"311\u00f79=34, 5\"),\n    @(\"261\u00f75=52, 1\", \"939\u00f75=187, 4\"),\n    @(\"577\u00f74=144, 1\", \"573\u00f75=114, 3\"),\n    @(\"903\u00f79=100, 3\", \"564\u00f78=70, 4\"),\n    @(\"545\u00f74=136, 1\", \"797\u00f75=159, 2\"),\n    @(\"576\u00f78=72, 0\", \"954\u00f73=318, 0\"),\n    @(\"744\u00f72=372, 0\", \"230\u00f76=38, 2\"),\n    @(\"793\u00f76=132, 1\", \"148\u00f76=24, 4\"),\n    @(\"528\u00f78=66, 0\", \"290\u00f75=58, 0\"),\n    @(\"108\u00f79=12, 0\", \"340\u00f74=85, 0\"),\n    @(\"749\u00f76=124, 5\", \"913\u00f76=152, 1\"),\n    @(\"224\u00f74=56, 0\", \"658\u00f78=82, 2\"),\n    @(\"260\u00f78=32, 4\", \"925\u00f75=185, 0\"),\n    @(\"607\u00f73=202, 1\", \"286\u00f76=47, 4\"),\n    @(\"927\u00f78=115, 7\", \"267\u00f72=133, 1\"),\n    @(\"923\u00f78=115, 3\", \"902\u00f79=100, 2\"),\n    @(\"350\u00f72=175, 0\", \"380\u00f78=47, 4\"),\n    @(\"121\u00f77=17, 2\", \"807\u00f79=89, 6\"),\n    @(\"219\u00f73=73, 0\", \"294\u00f73=98, 0\"),\n    @(\"991\u00f75=198, 1\", \"706\u00f72=353, 0\"),\n    @(\"683\u00f75=136, 3\", \"980\u00f73=326, 2\"),\n    @(\"254\u00f74=63, 2\", \"732\u00f78=91, 4\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $false, $true, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
